$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Resource Mix")

$ws.Range("E2").Value = "25 Lead,20 Sr Lead"
$ws.Range("E3").Value = "3 Sr Lead,2 Lead"
